# Apply cryptos-list update (commit: "Updated cryptos list on Fri Nov 22 20:14:52 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    # Leading apostrophe forces Excel to treat numeric/ambiguous-looking
    # text (e.g. "1.00", "0.401", "3.10") as literal text instead of
    # coercing it to a number. Resetting the Style afterwards clears the
    # resulting quote-prefix style so the cell keeps its original (default)
    # formatting, matching the source workbook.
    $ws.Range($addr).Value = "" + $text
    $ws.Range($addr).Style = "Normal"
}

Set-CellText "D2" "99.280.26"
Set-CellText "E2" "  +0.92%  "
Set-CellText "D3" "3.285.50"
Set-CellText "E3" "  -2.11%  "
Set-CellText "D5" "253.90"
Set-CellText "E5" "  -0.92%  "
Set-CellText "D6" "622.78"
Set-CellText "E6" "  -0.19%  "
Set-CellText "D7" "1.44"
Set-CellText "E7" "  +22.01%  "
Set-CellText "D8" "0.401"
Set-CellText "E8" "  +3.35%  "
Set-CellText "E9" "  -0.01%  "
Set-CellText "D10" "0.975"
Set-CellText "E10" "  +22.29%  "
Set-CellText "D11" "3.279.86"
Set-CellText "E11" "  -2.22%  "
Set-CellText "D12" "0.201"
Set-CellText "E12" "  +0.79%  "
Set-CellText "D13" "39.44"
Set-CellText "E13" "  +9.63%  "
Set-CellText "D14" "99.074.82"
Set-CellText "E14" "  +1.02%  "
Set-CellText "E15" "  +0.07%  "
Set-CellText "D16" "3.891.97"
Set-CellText "E16" "  -2.30%  "
Set-CellText "D17" "5.48"
Set-CellText "E17" "  -0.33%  "
Set-CellText "D18" "3.289.71"
Set-CellText "E18" "  -2.09%  "
Set-CellText "D19" "3.44"
Set-CellText "E19" "  -4.52%  "
Set-CellText "E20" "  +1.92%  "
Set-CellText "D21" "6.32"
Set-CellText "E21" "  +7.44%  "
Set-CellText "D22" "487.15"
Set-CellText "E22" "  -0.56%  "
Set-CellText "D23" "9.28"
Set-CellText "E23" "  +1.30%  "
Set-CellText "D24" "0.0000201"
Set-CellText "E24" "  -2.91%  "
Set-CellText "D25" "5.62"
Set-CellText "E25" "  -1.00%  "
Set-CellText "B26" "Stellar"
Set-CellText "C26" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText "D26" "0.336"
Set-CellText "E26" "  +37.67%  "
Set-CellText "B27" "Litecoin"
Set-CellText "C27" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-CellText "D27" "89.12"
Set-CellText "E27" "  +1.14%  "
Set-CellText "D28" "11.99"
Set-CellText "E28" "  -0.29%  "
Set-CellText "D29" "3.446.03"
Set-CellText "E29" "  -2.56%  "
Set-CellText "E30" "  -0.08%  "
Set-CellText "B31" "Hedera"
Set-CellText "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText "D31" "0.138"
Set-CellText "E31" "  +12.42%  "
Set-CellText "B32" "Cronos"
Set-CellText "C32" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText "D32" "0.190"
Set-CellText "E32" "  +1.94%  "
Set-CellText "D33" "10.35"
Set-CellText "E33" "  +11.92%  "
Set-CellText "D34" "0.999"
Set-CellText "E34" "  +0.03%  "
Set-CellText "D35" "28.01"
Set-CellText "E35" "  +2.83%  "
Set-CellText "B36" "PolygonEcosystemToken"
Set-CellText "C36" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-CellText "D36" "0.477"
Set-CellText "E36" "  +6.68%  "
Set-CellText "B37" "Kaspa"
Set-CellText "C37" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText "D37" "0.150"
Set-CellText "E37" "  -0.31%  "
Set-CellText "D38" "7.22"
Set-CellText "E38" "  -1.31%  "
Set-CellText "D39" "1.94"
Set-CellText "E39" "  -0.44%  "
Set-CellText "D40" "24.76"
Set-CellText "E40" "  -0.40%  "
Set-CellText "D41" "486.51"
Set-CellText "E41" "  -5.75%  "
Set-CellText "D42" "3.71"
Set-CellText "E42" "  +1.38%  "
Set-CellText "D43" "1.23"
Set-CellText "E43" "  -2.55%  "
Set-CellText "E44" "  -0.03%  "
Set-CellText "D45" "0.772"
Set-CellText "E45" "  -0.97%  "
Set-CellText "D46" "3.10"
Set-CellText "E46" "  -5.54%  "
Set-CellText "D47" "1.95"
Set-CellText "E47" "  +1.16%  "
Set-CellText "D48" "157.71"
Set-CellText "E48" "  -1.38%  "
Set-CellText "D49" "0.848"
Set-CellText "E49" "  +7.27%  "
Set-CellText "D50" "7.29"
Set-CellText "E50" "  +15.23%  "
Set-CellText "D51" "4.71"
Set-CellText "E51" "  +4.54%  "
